# Update the RP shield BOM: row 13 (item #7, RP-E1/RP-E2 connector) gets a
# real manufacturer + part number, and its description is clarified to call
# out that it's a "raspberry pi style stacking header".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

$ws.Range("D13").Value2 = "Adafruit Industries LLC"
$ws.Range("E13").Value2 = "1979  [1528-1783-ND (DigiKey)]"
$ws.Range("F13").Value2 = "2x13 (26pin) RP connector (raspberry pi style stacking header, long) BOTTOM MOUNT!"

# Leave the view/selection the way the author left it after editing: scrolled
# back to the top of the sheet with E13 (the last-touched cell) selected.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select() | Out-Null
